$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23 corresponds to "urban_nr" - incorporate new urban net returns data.
# Clear the older (1997, 2002) values in E23:F23
$ws.Range("E23").Value = $null
$ws.Range("F23").Value = $null

# Update 2007, 2012 values with the newly computed figures
$ws.Range("G23").Value = 10400.7021484375
$ws.Range("H23").Value = 9830.865234375

# Add the new 2015 value
$ws.Range("I23").Value = 8696.31640625
